# Fixes error in ggplot2 cheatsheet
#
# Corrects a handful of typos/smart-quote artifacts in the code-sample
# text boxes on slide 1, and fixes a copy/paste mistake where the
# "discrete x, continuous y" example reused the "e" plot object /
# continuous-x example instead of introducing "f" with a discrete x.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# 1) "a + geom_path(...)" box: curly closing quote -> straight quote
$shp138 = Get-ShapeById $s 138
$tr138 = $shp138.TextFrame.TextRange
for ($j = 1; $j -le $tr138.Paragraphs().Count; $j++) {
    $para = $tr138.Paragraphs($j)
    if ($para.Text -like "*geom_path*") {
        $para.Runs(2).Text = 'lineend="butt", linejoin="round", linemitre=1'
    }
}

# 2) "e + geom_rug(...)" / "e <- ggplot(...)" box: curly quotes -> straight quotes
$shp228 = Get-ShapeById $s 228
$tr228 = $shp228.TextFrame.TextRange
for ($j = 1; $j -le $tr228.Paragraphs().Count; $j++) {
    $para = $tr228.Paragraphs($j)
    if ($para.Text -like "*geom_rug*") {
        $para.Runs(2).Text = 'sides = "bl"'
    }
}

# 3) "discrete x , continuous y" box: was a copy of the continuous/continuous
#    example ("e <- ggplot(mpg, aes(cty, hwy))"); should introduce "f" with
#    a discrete x-axis variable ("class"). This textbox auto-sizes to fit
#    its text (spAutoFit), so temporarily disable that while editing the
#    run to avoid PowerPoint silently nudging the shape's stored height as
#    a side effect of the text edit.
$shp229 = Get-ShapeById $s 229
$origAutoSize229 = $shp229.TextFrame.AutoSize
$origHeight229 = $shp229.Height
$shp229.TextFrame.AutoSize = 0
$tr229 = $shp229.TextFrame.TextRange
for ($j = 1; $j -le $tr229.Paragraphs().Count; $j++) {
    $para = $tr229.Paragraphs($j)
    if ($para.Text -like "*ggplot*") {
        $para.Runs(1).Text = "f <- ggplot(mpg, aes(class, hwy))"
    }
}
$shp229.TextFrame.AutoSize = $origAutoSize229
$shp229.Height = $origHeight229 + 0.00001

# 4) "f + geom_dotplot(...)" and "f + geom_violin(...)" box: curly quotes -> straight quotes
$shp230 = Get-ShapeById $s 230
$tr230 = $shp230.TextFrame.TextRange
for ($j = 1; $j -le $tr230.Paragraphs().Count; $j++) {
    $para = $tr230.Paragraphs($j)
    if ($para.Text -like "*geom_dotplot*") {
        $para.Runs(2).Text = 'binaxis = "y", stackdir = "center"'
    }
    if ($para.Text -like "*geom_violin*") {
        $para.Runs(2).Text = 'scale = "area"'
    }
}
